# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.539.68'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.36'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -1.40%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.46'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4632'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3846'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.02'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07907'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('E10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9932'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('E11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.45'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.847.12'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.912'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.101'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.84'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06654'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('E18').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.05'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.28%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '27.539.92'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.379'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.303'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.60%  '
$ws.Range('E25').ClearFormats()
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B26').ClearFormats()
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C26').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.062.98'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('E26').ClearFormats()
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Monero'
$ws.Range('B27').ClearFormats()
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C27').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.98'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('E27').ClearFormats()
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('B28').ClearFormats()
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.47'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.84%  '
$ws.Range('E28').ClearFormats()
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('B29').ClearFormats()
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C29').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.099'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E29').ClearFormats()
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('B30').ClearFormats()
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C30').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.394'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('E30').ClearFormats()
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('B31').ClearFormats()
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C31').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.63'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('E31').ClearFormats()
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('B32').ClearFormats()
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C32').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9747'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('E32').ClearFormats()
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Stellar'
$ws.Range('B33').ClearFormats()
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C33').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09406'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.76%  '
$ws.Range('E33').ClearFormats()
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('B34').ClearFormats()
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C34').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.574'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.03%  '
$ws.Range('E34').ClearFormats()
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('B35').ClearFormats()
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C35').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.278'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E35').ClearFormats()
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('B36').ClearFormats()
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C36').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.344'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('E36').ClearFormats()
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Hedera'
$ws.Range('B37').ClearFormats()
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C37').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06011'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('E37').ClearFormats()
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('B38').ClearFormats()
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C38').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02220'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('E38').ClearFormats()
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('B39').ClearFormats()
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C39').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.278'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('E39').ClearFormats()
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('B40').ClearFormats()
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C40').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.178'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('E40').ClearFormats()
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5889'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('E41').ClearFormats()
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Algorand'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1859'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('E42').ClearFormats()
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Aptos'
$ws.Range('B43').ClearFormats()
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C43').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('E43').ClearFormats()
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('B44').ClearFormats()
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C44').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.253'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5578'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.29%  '
$ws.Range('E45').ClearFormats()
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('B46').ClearFormats()
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C46').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.20'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('E46').ClearFormats()
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('B47').ClearFormats()
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C47').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.896'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('E47').ClearFormats()
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Cronos'
$ws.Range('B48').ClearFormats()
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C48').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06687'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('E48').ClearFormats()
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Quant'
$ws.Range('B49').ClearFormats()
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C49').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '110.59'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('E49').ClearFormats()
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'EOS'
$ws.Range('B50').ClearFormats()
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('C50').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.052'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('E50').ClearFormats()
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('B51').ClearFormats()
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C51').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.47%  '
$ws.Range('E51').ClearFormats()
